$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# --- Row 27 & 29: role cell correction (System Analyst -> Software Architect) ---
$ws.Range("B27").Value = "Software Architect"
$ws.Range("C28").Value = 43893
$ws.Range("B29").Value = "Software Architect"
$ws.Range("C29").Value = 43893

# --- New activity rows 31..36 (filled first so the shared-string table gets the
#     same append order as the authored workbook: strings 71-76 come from column
#     A of rows 31-36, 77 ("15 min") from F32, and 78 last from A30) ---

# Row 31
$ws.Range("A31").Value = "Lav SD0804 angivStraksafskrivning"
$ws.Range("B31").Value = "Software Architect"
$ws.Range("C31").Value = 43894
$ws.Range("D31").Value = 0.52083333333333337
$ws.Range("E31").Value = 0.59375

# Row 32
$ws.Range("A32").Value = "Indsaml data til UC09 Beregn resultat uden renter"
$ws.Range("B32").Value = "Requirement Specifier"
$ws.Range("C32").Value = 43894
$ws.Range("D32").Value = 0.59375
$ws.Range("E32").Value = 0.625

# Row 33
$ws.Range("A33").Value = "Lav UC09 Beregn resultat før renter"
$ws.Range("B33").Value = "Requirement Specifier"
$ws.Range("C33").Value = 43894
$ws.Range("D33").Value = 0.625
$ws.Range("E33").Value = 0.63541666666666663

# Row 34
$ws.Range("A34").Value = "Lav DOM09 Beregn resultat før renter"
$ws.Range("B34").Value = "Requirement Specifier"
$ws.Range("C34").Value = 43894
$ws.Range("D34").Value = 0.63541666666666663
$ws.Range("E34").Value = 0.64583333333333337

# Row 35
$ws.Range("A35").Value = "Lav ADT09a Beregn resultat før renter"
$ws.Range("B35").Value = "Requirement Specifier"
$ws.Range("C35").Value = 43894
$ws.Range("D35").Value = 0.64583333333333337
$ws.Range("E35").Value = 0.65625

# Row 36
$ws.Range("A36").Value = "Review UC10, DOM10 og ADT10a"
$ws.Range("B36").Value = "Reviewer"
$ws.Range("C36").Value = 43894
$ws.Range("D36").Value = 0.65625
$ws.Range("E36").Value = 0.67361111111111116

# F column "15 min" notes for rows 32-36 (first use on F32 mints the shared string)
$ws.Range("F32").Value = "15 min"
$ws.Range("F33").Value = "15 min"
$ws.Range("F34").Value = "15 min"
$ws.Range("F35").Value = "15 min"
$ws.Range("F36").Value = "15 min"

# Match the original cell alignment style used for the F column: rows 30-32 use
# the "blank" xf 18 (horizontal+vertical center) already seen elsewhere in the
# column, rows 33-36 use xf 20 (horizontal-center only).
$ws.Range("F7").Copy()
$ws.Range("F30:F32").PasteSpecial(-4122)
$ws.Range("F5").Copy()
$ws.Range("F33:F36").PasteSpecial(-4122)

# Row 30 (filled last so its new string becomes the final new shared-string entry)
$ws.Range("A30").Value = "Lav SD0802 angivLineaerAfskrivning"
$ws.Range("B30").Value = "Software Architect"
$ws.Range("C30").Value = 43894
$ws.Range("D30").Value = 0.35416666666666669
$ws.Range("E30").Value = 0.39583333333333331
# H30 used to carry the running-total formula; the author cleared it out for this row.
$ws.Range("H30").ClearContents()

# --- Row 40: extend the G/H running-total formulas down one more row (mirrors
#     the shared formula range G9:G39 -> G9:G40 and SUM(G$5:G39) -> SUM(G$5:G40)) ---
$ws.Range("G40").Formula = "=E40-D40"
$ws.Range("G39").Copy()
$ws.Range("G40").PasteSpecial(-4122)

$ws.Range("H40").Formula = "=SUM(G`$5:G40)"
$ws.Range("H39").Copy()
$ws.Range("H40").PasteSpecial(-4122)

$ws.Rows.Item(40).RowHeight = 19.8

# --- Row 55: new trailing blank row matching the style of the rows above it ---
$ws.Range("C54").Copy()
$ws.Range("C55").PasteSpecial(-4122)

# --- View state: selection moved to B23 (scrolled down toward that area) ---
$ws.Activate()
$ws.Range("B23").Select()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
